$p = $ppt.ActivePresentation

# --- 1. Remove the old slide 8 (sldId 262) from the slide show ---
# This is the simple "Title + big picture" slide that duplicated the
# "Question 4" content; the slide that follows it (sldId 265, the detailed
# "Question 4 - Birth Rate and Happiness Score" slide) shifts up to
# become the new slide 8.
$p.Slides.Item(8).Delete()

# The former slide 9 (sldId 265) is now at index 8.
$s = $p.Slides.Item(8)

# --- 2. Update "Rectangle 3" Pearson r textbox ---
$rect = $s.Shapes.Item(2)
$rect.TextFrame.TextRange.Text = "Pearson r correlation coefficient = 0.73"

# --- 3. Update "TextBox 6" - change only the bold/red middle run's text,
# keeping the surrounding plain-text runs and the bold/red formatting
# of the replaced run intact. ---
$tb6 = $s.Shapes.Item(3)
$tb6tr = $tb6.TextFrame.TextRange
$oldMiddle = "a weak-medium negative relationship "
$fullText = $tb6tr.Text
$startPos = $fullText.IndexOf($oldMiddle) + 1
$midRange = $tb6tr.Characters($startPos, $oldMiddle.Length)
$midRange.Text = "a strong positive relationship "

# --- 4. Update "TextBox 9" - Regression Analysis numbers: reorder /
# replace the paragraphs and append two new ones. ---
$tb9 = $s.Shapes.Item(7)
$tb9tr = $tb9.TextFrame.TextRange
$tb9tr.Text = "Observations  137`rR Square 0.530437094  `rP-value  6.46E-24  `rCoefficients  0.0000437763`rEach `$10K = 0.437763`rEach `$30K = 1.313289 "
# Italicize "P-value" (paragraph 3, starts at char 42, length 7)
$tb9tr.Characters(42, 7).Font.Italic = $true
# Italicize "Coefficients" (paragraph 4, starts at char 62, length 12)
$tb9tr.Characters(62, 12).Font.Italic = $true

# --- 5. Move the picture (Content Placeholder 4) to resize/reposition it
# and bring it to the front of the z-order (so it is the last shape in
# the slide's shape tree, drawn after the text boxes). ---
$pic = $s.Shapes.Item(1)
$pic.Left = 246064 / 12700.0
$pic.Width = 7340065 / 12700.0
$pic.Height = 5505050 / 12700.0
$pic.ZOrder(0)
